$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1446.5454
$ws.Range("I80").Value = 640.25
$ws.Range("J80").Value = 3596.6667
$ws.Range("K80").Value = 1920.75
$ws.Range("L80").Value = 10790.0001
$ws.Range("M80").Value = -922.75
$ws.Range("N80").Value = -12786.0001
$ws.Range("H83").Value = 1446.5454
$ws.Range("I83").Value = 640.25
$ws.Range("J83").Value = 3596.6667
$ws.Range("K83").Value = 5762.25
$ws.Range("L83").Value = 32370.0003
$ws.Range("M83").Value = -770.25
$ws.Range("N83").Value = -42354.0003
$ws.Range("H133").Value = 46577.777
$ws.Range("J133").Value = 46577.777
$ws.Range("L133").Value = 46577.777
$ws.Range("N133").Value = -56697.777
$ws.Range("H138").Value = 6669293.5
$ws.Range("J138").Value = 10002979
$ws.Range("L138").Value = 30008937
$ws.Range("N138").Value = -30019217

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18272.404
$ws.Range("I32").Value = 4253.507
$ws.Range("J32").Value = 188835.67
$ws.Range("K32").Value = 4253.507
$ws.Range("L32").Value = 188835.67
$ws.Range("M32").Value = -3966.507
$ws.Range("N32").Value = -189409.67
$ws.Range("H61").Value = 1768.7037
$ws.Range("I61").Value = 1464.7916
$ws.Range("J61").Value = 4200
$ws.Range("K61").Value = 1464.7916
$ws.Range("L61").Value = 4200
$ws.Range("M61").Value = -1252.7916
$ws.Range("N61").Value = -4624
$ws.Range("H122").Value = 2145
$ws.Range("I122").Value = 1806.25
$ws.Range("K122").Value = 5418.75
$ws.Range("M122").Value = -2968.75
$ws.Range("H132").Value = 2088.1228
$ws.Range("I132").Value = 1878.28
$ws.Range("J132").Value = 3587
$ws.Range("K132").Value = 5634.84
$ws.Range("L132").Value = 10761
$ws.Range("M132").Value = -3104.84
$ws.Range("N132").Value = -15821
$ws.Range("H133").Value = 43753.668
$ws.Range("J133").Value = 43753.668
$ws.Range("L133").Value = 43753.668
$ws.Range("N133").Value = -48813.668
$ws.Range("H136").Value = 1768.7037
$ws.Range("I136").Value = 1464.7916
$ws.Range("J136").Value = 4200
$ws.Range("K136").Value = 4394.3748
$ws.Range("L136").Value = 12600
$ws.Range("M136").Value = -1844.3748
$ws.Range("N136").Value = -17700
$ws.Range("H139").Value = 56153.75
$ws.Range("J139").Value = 56153.75
$ws.Range("L139").Value = 56153.75
$ws.Range("N139").Value = -66433.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 10723.333
$ws.Range("I5").Value = 12680
$ws.Range("K5").Value = 12680
$ws.Range("M5").Value = -12567
$ws.Range("H22").Value = 254.4
$ws.Range("I22").Value = 238.22223
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 238.22223
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -65.22223
$ws.Range("N22").Value = -746
$ws.Range("H59").Value = 58450
$ws.Range("J59").Value = 58450
$ws.Range("L59").Value = 58450
$ws.Range("N59").Value = -60144

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 63.5
$ws.Range("I7").Value = 50.333332
$ws.Range("J7").Value = 76.666664
$ws.Range("K7").Value = 50.333332
$ws.Range("L7").Value = 76.666664
$ws.Range("M7").Value = 62.666668
$ws.Range("N7").Value = -302.666664
$ws.Range("H99").Value = 27780440
$ws.Range("I99").Value = 2616.25
$ws.Range("K99").Value = 2616.25
$ws.Range("M99").Value = -1118.25
$ws.Range("H107").Value = 429.9091
$ws.Range("I107").Value = 289
$ws.Range("J107").Value = 676.5
$ws.Range("K107").Value = 289
$ws.Range("L107").Value = 676.5
$ws.Range("M107").Value = 1631
$ws.Range("N107").Value = -4516.5
$ws.Range("H126").Value = 27780440
$ws.Range("I126").Value = 2616.25
$ws.Range("K126").Value = 7848.75
$ws.Range("M126").Value = -5378.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 3895.238
$ws.Range("I69").Value = 1000
$ws.Range("J69").Value = 4200
$ws.Range("K69").Value = 3000
$ws.Range("L69").Value = 12600
$ws.Range("M69").Value = -2189
$ws.Range("N69").Value = -14222
$ws.Range("H72").Value = 3895.238
$ws.Range("I72").Value = 1000
$ws.Range("J72").Value = 4200
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 37800
$ws.Range("M72").Value = -4944
$ws.Range("N72").Value = -45912
$ws.Range("H105").Value = 5000
$ws.Range("J105").Value = 5000
$ws.Range("L105").Value = 15000
$ws.Range("N105").Value = -20242

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 11328.75
$ws.Range("J123").Value = 11328.75
$ws.Range("L123").Value = 11328.75
$ws.Range("N123").Value = -16228.75
$ws.Range("H126").Value = 2923.375
$ws.Range("I126").Value = 2740.1428
$ws.Range("J126").Value = 2998.8235
$ws.Range("K126").Value = 8220.428400000001
$ws.Range("L126").Value = 8996.470499999999
$ws.Range("M126").Value = -5750.428400000001
$ws.Range("N126").Value = -13936.4705
$ws.Range("H132").Value = 3193.5476
$ws.Range("I132").Value = 2779.4119
$ws.Range("J132").Value = 4953.625
$ws.Range("K132").Value = 8338.235700000001
$ws.Range("L132").Value = 14860.875
$ws.Range("M132").Value = -5808.235700000001
$ws.Range("N132").Value = -19920.875
$ws.Range("H138").Value = 69602.89999999999
$ws.Range("J138").Value = 69602.89999999999
$ws.Range("L138").Value = 69602.89999999999
$ws.Range("N138").Value = -79882.89999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3793.5293
$ws.Range("I7").Value = 6990
$ws.Range("J7").Value = 3593.75
$ws.Range("K7").Value = 6990
$ws.Range("L7").Value = 3593.75
$ws.Range("M7").Value = -6878
$ws.Range("N7").Value = -3817.75
$ws.Range("H40").Value = 3666.6667
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3666.6667
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3666.6667
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -3938.6667
$ws.Range("H126").Value = 3793.5293
$ws.Range("I126").Value = 6990
$ws.Range("J126").Value = 3593.75
$ws.Range("K126").Value = 20970
$ws.Range("L126").Value = 10781.25
$ws.Range("M126").Value = -18500
$ws.Range("N126").Value = -15721.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 9166.666999999999
$ws.Range("J22").Value = 9166.666999999999
$ws.Range("L22").Value = 9166.666999999999
$ws.Range("N22").Value = -9752.666999999999
$ws.Range("H57").Value = 42863.332
$ws.Range("I57").Value = 55000
$ws.Range("J57").Value = 36795
$ws.Range("K57").Value = 55000
$ws.Range("L57").Value = 36795
$ws.Range("M57").Value = -54246
$ws.Range("N57").Value = -38303
$ws.Range("H122").Value = 2726
$ws.Range("I122").Value = 2452
$ws.Range("K122").Value = 7356
$ws.Range("M122").Value = -4906
$ws.Range("H126").Value = 77910.69500000001
$ws.Range("I126").Value = 143199.28
$ws.Range("K126").Value = 429597.84
$ws.Range("M126").Value = -427127.84
